$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Orders" ---
$ws = $wb.Worksheets.Item("All Orders")

# Insert a new row above row 2 - shifts existing order rows (old row2->3, old row3->4)
$ws.Rows.Item(2).Insert()

# Fill in the new order row (row 2) - a new order placed at 10:38
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = "2026-01-13 10:38"
$ws.Cells.Item(2, 3).Value = "Pooja"
# Phone number must stay text (not be auto-converted to a number)
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "9096648553"
$ws.Cells.Item(2, 5).Value = "A 1608, Pune 411045"
$ws.Cells.Item(2, 6).Value = "Girl Haldi Kunku Set x1"
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""

# --- Sheet 2: "Daily Summary" ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(2, 2).Value = 3
